$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.426.40"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.013.65"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "262.80"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.76%  "
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "56.11"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -6.69%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").Value = "  -2.10%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.33"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("D13").Value = "2.310.45"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E14").Value = "  -4.79%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.04"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -7.16%  "
$ws.Range("E16").Value = "  -3.59%  "
$ws.Range("D17").Value = "2.023.08"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "37.295.92"
$ws.Range("E18").Value = "  +0.43%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "69.75"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("E21").Value = "  -0.13%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "228.46"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.69"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +8.82%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -0.69%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "164.72"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.43%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.01"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.61%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.70"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.128"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -10.32%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -3.18%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0650"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  +0.00%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("E37").Value = "  -0.22%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.34"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.19%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.20"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.22%  "
$ws.Range("E40").Value = "  +4.46%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.22"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("E42").Value = "  -3.66%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0213"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "1.395.72"
$ws.Range("E44").Value = "  +1.04%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "90.25"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -5.83%  "
$ws.Range("E47").Value = "  -1.83%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.08"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").Value = "2.202.06"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  -3.92%  "
